# Insert a new data row for the "RDU" (Durham, United States) colo just
# above the existing "ADL" row (row 299), shifting the Oceania/Pacific
# block (ADL..PPT) down by one row, matching the updated generated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 299 (shifts rows 299:310 down to 300:311).
$ws.Rows.Item(299).Insert()

# Populate the new row with the RDU colo data.
$ws.Range("A299").Value = "RDU"
$ws.Range("B299").Value = "Durham, United States"
$ws.Range("C299").Value = "RDU"
$ws.Range("D299").Value = 35.93543
$ws.Range("E299").Value = -78.88075000000001
$ws.Range("F299").Value = "US"
$ws.Range("G299").Value = "North America"
$ws.Range("H299").Value = "Durham"

# The inserted row picks up formatting from the row above by default, which
# lacks the bold/bordered style used by column A throughout the table.
# Copy that style from the (now shifted) following row so A299 matches the
# rest of the "colo" column.
$ws.Range("A300").Copy()
$ws.Range("A299").PasteSpecial(-4122)
